# lab_2_paperweights.xlsx - "new data time stamps"
# The lab data table (A1:B119, header in row 1: area_cm2 / mass_g) was
# re-sorted in ascending order by column A ("area_cm2") using Excel's
# Data > Sort feature on a whole-column selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mimic selecting the whole A:B columns (as one would before invoking
# Data > Sort from the ribbon) so the sheet's recorded selection matches
# the post-sort state.
$ws.Columns("A:B").Select()

# Sort the data range (including the header row) ascending by column A.
# Header:=1 (xlYes) keeps row 1 ("area_cm2"/"mass_g") pinned in place.
$dataRange = $ws.Range("A1:B119")
$key1 = $ws.Range("A1")
$dataRange.Sort($key1, 1, $null, $null, 1, $null, 1, 1)
